$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the validation script cell (K7) to reference the new dynamic property
# plugin script instead of the inline "Hello world" script.
$ws.Range("K7").Value = "test-dynamic.py"

# Move the active selection to the edited cell, matching the saved workbook state.
$ws.Range("K7").Select()
